$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 12
$ws.Range("AB2").Value = 9.5
$ws.Range("AH2").Value = 11
$ws.Range("AM2").Value = 41

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("I5").Value = 4.1
$ws.Range("L5").Value = 4.5
$ws.Range("S5").Value = 3.85
$ws.Range("AD5").Value = 19
$ws.Range("AE5").Value = 21
$ws.Range("AI5").Value = 17

# Row 10
$ws.Range("Q10").Value = 2.4
$ws.Range("R10").Value = 1.53
$ws.Range("S10").Value = 3.8
$ws.Range("T10").Value = 1.26
$ws.Range("AR10").Value = 1.85
$ws.Range("AS10").Value = 2

# Row 11
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 1.95
$ws.Range("L11").Value = 4.33
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 7.5
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 3.75
$ws.Range("U11").Value = 5
$ws.Range("V11").Value = 1.17
$ws.Range("W11").Value = 1.53
$ws.Range("X11").Value = 2.38
$ws.Range("Y11").Value = 2.1
$ws.Range("Z11").Value = 1.67
$ws.Range("AA11").Value = 6
$ws.Range("AB11").Value = 9
$ws.Range("AC11").Value = 10
$ws.Range("AD11").Value = 19
$ws.Range("AF11").Value = 41
$ws.Range("AG11").Value = 7.5
$ws.Range("AH11").Value = 7
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 81
$ws.Range("AL11").Value = 7.5
$ws.Range("AM11").Value = 15
$ws.Range("AN11").Value = 12
$ws.Range("AO11").Value = 41
$ws.Range("AP11").Value = 34
$ws.Range("AQ11").Value = 41
$ws.Range("AR11").Value = 1.85
$ws.Range("AS11").Value = 2

# Row 12
$ws.Range("G12").Value = 1.12
$ws.Range("H12").Value = 7.5
$ws.Range("I12").Value = 25
$ws.Range("J12").Value = 1.42
$ws.Range("K12").Value = 2.92
$ws.Range("L12").Value = 15
$ws.Range("N12").Value = 10.75
$ws.Range("O12").Value = 1.14
$ws.Range("P12").Value = 5.2
$ws.Range("Q12").Value = 1.45
$ws.Range("R12").Value = 2.65
$ws.Range("U12").Value = 2.07
$ws.Range("V12").Value = 1.7
$ws.Range("W12").Value = 1.25
$ws.Range("X12").Value = 3.7
$ws.Range("AA12").Value = 6.7
$ws.Range("AB12").Value = 6
$ws.Range("AC12").Value = 12.5
$ws.Range("AD12").Value = 6.1
$ws.Range("AG12").Value = 10.75
$ws.Range("AH12").Value = 18.5
$ws.Range("AI12").Value = 50
$ws.Range("AJ12").Value = 300
$ws.Range("AL12").Value = 45
$ws.Range("AM12").Value = 350
$ws.Range("AN12").Value = 90
$ws.Range("AP12").Value = 600
$ws.Range("AQ12").Value = 350

# Row 17
$ws.Range("G17").Value = 2.1
$ws.Range("I17").Value = 2.9
$ws.Range("L17").Value = 3.4
$ws.Range("N17").Value = 15
$ws.Range("S17").Value = 1.93
$ws.Range("T17").Value = 1.88
$ws.Range("AA17").Value = 11
$ws.Range("AC17").Value = 9
$ws.Range("AD17").Value = 21
$ws.Range("AE17").Value = 15
$ws.Range("AL17").Value = 15
$ws.Range("AM17").Value = 19
$ws.Range("AN17").Value = 12
$ws.Range("AO17").Value = 34

# Row 18
$ws.Range("G18").Value = 2.2
$ws.Range("I18").Value = 3.1
$ws.Range("J18").Value = 2.77
$ws.Range("L18").Value = 3.6
$ws.Range("AA18").Value = 9
$ws.Range("AB18").Value = 12
$ws.Range("AD18").Value = 23
$ws.Range("AE18").Value = 16.5
$ws.Range("AI18").Value = 12
$ws.Range("AL18").Value = 10.75
$ws.Range("AM18").Value = 17.5
$ws.Range("AN18").Value = 10.75
$ws.Range("AO18").Value = 40
$ws.Range("AP18").Value = 25
$ws.Range("AQ18").Value = 29

# Row 19
$ws.Range("G19").Value = 3.6
$ws.Range("H19").Value = 3.45
$ws.Range("I19").Value = 1.91
$ws.Range("J19").Value = 4
$ws.Range("K19").Value = 2.12
$ws.Range("L19").Value = 2.52
$ws.Range("O19").Value = 1.24
$ws.Range("P19").Value = 3.3
$ws.Range("Q19").Value = 1.72
$ws.Range("R19").Value = 1.88
$ws.Range("U19").Value = 2.67
$ws.Range("V19").Value = 1.36
$ws.Range("Y19").Value = 1.62
$ws.Range("Z19").Value = 2.02
$ws.Range("AA19").Value = 12
$ws.Range("AB19").Value = 21
$ws.Range("AC19").Value = 12
$ws.Range("AD19").Value = 55
$ws.Range("AF19").Value = 32
$ws.Range("AG19").Value = 11.25
$ws.Range("AH19").Value = 6.7
$ws.Range("AK19").Value = 350
$ws.Range("AM19").Value = 9.75
$ws.Range("AN19").Value = 8.25
$ws.Range("AO19").Value = 17
$ws.Range("AP19").Value = 14.5
$ws.Range("AQ19").Value = 24

# Row 27
$ws.Range("AA27").Value = 7.9
$ws.Range("AB27").Value = 9.5
$ws.Range("AI27").Value = 14

# Row 28
$ws.Range("Q28").Value = 1.98
$ws.Range("R28").Value = 1.88
$ws.Range("U28").Value = 3.4
$ws.Range("V28").Value = 1.3

# Row 29
$ws.Range("U29").Value = 2
$ws.Range("V29").Value = 1.73
$ws.Range("AA29").Value = 26

# Row 30
$ws.Range("G30").Value = 1.95
$ws.Range("H30").Value = 3.1
$ws.Range("I30").Value = 3.9
$ws.Range("J30").Value = 2.52
$ws.Range("K30").Value = 2
$ws.Range("L30").Value = 4.45
$ws.Range("O30").Value = 1.39
$ws.Range("P30").Value = 2.57
$ws.Range("Q30").Value = 2.12
$ws.Range("R30").Value = 1.57
$ws.Range("U30").Value = 3.5
$ws.Range("V30").Value = 1.21
$ws.Range("W30").Value = 1.44
$ws.Range("X30").Value = 2.42
$ws.Range("Y30").Value = 1.9
$ws.Range("Z30").Value = 1.72
$ws.Range("AA30").Value = 6.2
$ws.Range("AB30").Value = 8.75
$ws.Range("AC30").Value = 8.5
$ws.Range("AD30").Value = 17.5
$ws.Range("AE30").Value = 17
$ws.Range("AF30").Value = 32
$ws.Range("AG30").Value = 7.5
$ws.Range("AH30").Value = 6.1
$ws.Range("AI30").Value = 16.5
$ws.Range("AJ30").Value = 90
$ws.Range("AK30").Value = 900
$ws.Range("AL30").Value = 9
$ws.Range("AM30").Value = 20
$ws.Range("AN30").Value = 13.5
$ws.Range("AO30").Value = 65
$ws.Range("AP30").Value = 45
$ws.Range("AQ30").Value = 55
